$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("D4").Value = "Circuitos Elétricos 2"
$ws.Range("E5").Value = "EAP"
$ws.Range("C7").Value = "-"
$ws.Range("F8").Value = "-"
$ws.Range("C9").Value = "EAP"
$ws.Range("F9").Value = "Circuitos Elétricos 2"
